# Commit: "Sat, May 16, 2020  9:05:51 PM"
#
# The underlying OOXML diff swaps the contents of ppt/theme/theme1.xml
# (the Slide Master's theme -- originally the "Integral" color scheme)
# and ppt/theme/theme2.xml (the Notes Master's theme -- originally the
# stock "Office Theme" color scheme): after the edit, theme1.xml carries
# the "Office" palette and theme2.xml carries the "Integral" palette
# (the font scheme and format scheme are identical between the two
# themes, so only the 12 color-scheme slots actually change).
#
# The PowerPoint object model only exposes a writable Theme off the
# Slide Master (SlideMaster.Theme / NotesMaster.Theme / HandoutMaster.Theme
# all resolve to the same single presentation theme), so we reproduce the
# reachable half of that swap: repainting the presentation's theme color
# scheme from the "Integral" palette to the stock "Office" palette via
# ThemeColorScheme, using the exact RGB values from the target theme.

$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$tcs = $sm.Theme.ThemeColorScheme

# ThemeColorScheme.Item index order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
